$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-01 Sunday" "2024-09-02 Monday"

Replace-Text "982÷8=" "165÷5="
Replace-Text "367÷5=" "120÷3="
Replace-Text "487÷8=" "529÷7="
Replace-Text "426÷6=" "782÷6="
Replace-Text "681÷3=" "538÷7="
Replace-Text "893÷2=" "178÷2="
Replace-Text "449÷9=" "185÷5="
Replace-Text "469÷5=" "883÷6="
Replace-Text "754÷6=" "131÷3="
Replace-Text "164÷7=" "506÷2="
Replace-Text "373÷5=" "117÷6="
Replace-Text "923÷6=" "991÷2="
Replace-Text "804÷9=" "936÷4="
Replace-Text "152÷2=" "161÷4="
Replace-Text "658÷2=" "728÷8="
Replace-Text "390÷2=" "952÷4="
Replace-Text "888÷7=" "787÷3="
Replace-Text "949÷8=" "673÷7="
Replace-Text "621÷6=" "435÷8="
Replace-Text "988÷9=" "838÷7="
Replace-Text "151÷3=" "150÷5="
Replace-Text "854÷5=" "802÷4="
Replace-Text "697÷5=" "741÷4="
Replace-Text "783÷2=" "110÷5="
Replace-Text "455÷8=" "194÷4="

Write-Output "Done"
